$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G2").Value = "2016-09-02 18:17:33"
$wsOverview.Range("G3").Value = "2016-09-02 18:17:33"

$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H2").Value = "2016-09-02 18:17:27"
$wsZhCn.Range("H3").Value = "2016-09-02 18:17:27"
$wsZhCn.Range("K2").Value = "2016-09-02 18:17:43"
$wsZhCn.Range("K3").Value = "2016-09-02 18:17:43"

$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("H2").Value = "2016-09-02 18:17:33"
$wsDeDe.Range("H3").Value = "2016-09-02 18:17:33"
$wsDeDe.Range("K2").Value = "2016-09-02 18:17:52"
$wsDeDe.Range("K3").Value = "2016-09-02 18:17:52"
